# Atualização de bases das ligas, do dia: 12-04-2024 às 20:28
# The source data rows got reshuffled (their unique match ids moved to a
# different row) while the running index in column A stayed put. For every
# affected row-group below we swap/rotate the B:AC payload between rows,
# keeping column A (and C/D/E, which are identical across the group anyway)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2) {
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Simple pairwise swaps (B:AC payload exchanged between the two rows)
Swap-Rows $ws 27 28
Swap-Rows $ws 47 48
Swap-Rows $ws 104 105
Swap-Rows $ws 143 145
Swap-Rows $ws 153 154
Swap-Rows $ws 214 215

# Three-way rotation: 148 <- 150, 149 <- 148, 150 <- 149
$range148 = $ws.Range("B148:AC148")
$range149 = $ws.Range("B149:AC149")
$range150 = $ws.Range("B150:AC150")

$v148 = $range148.Value()
$v149 = $range149.Value()
$v150 = $range150.Value()

$range149.Value = $v148
$range150.Value = $v149
$range148.Value = $v150
